# Generate Report for Archive
#
# The localization-status report is regenerated: the two outstanding
# files have moved on from "Ready for handoff" into active translation,
# so every cell carrying the old status text is refreshed, and the
# status/locale columns -- originally sized to fit the longer
# "Ready for handoff" string -- are re-fit to the shorter replacement.

$wb = $excel.ActiveWorkbook

$oldStatus = "Ready for handoff"
$newStatus = "In Translation"

# Replace the status text wherever it appears (Overview!E2:F3 and the
# "Status" column on each per-locale detail sheet). Put the string
# literal on the left of -eq so PowerShell doesn't coerce it to a
# boolean when a cell (e.g. a TRUE/FALSE flag column) holds a bool.
foreach ($ws in $wb.Worksheets) {
    $used = $ws.UsedRange
    foreach ($cell in $used.Cells) {
        if ($oldStatus -eq $cell.Value2) {
            $cell.Value = $newStatus
        }
    }
}

# Re-fit the columns that held the status text -- the new text is
# shorter, so the columns shrink to match.
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Columns.Item(5).ColumnWidth = 12.5   # column E ("zh-cn")
$wsOverview.Columns.Item(6).ColumnWidth = 12.5   # column F ("de-de")

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Columns.Item(3).ColumnWidth = 12.5       # column C ("Status")

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Columns.Item(3).ColumnWidth = 12.5       # column C ("Status")
